$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update individual cell values to the new curated dimensions
$ws.Range("D2").Value = "sdmx-dimension:refArea"
$ws.Range("F2").Value = "iaest-measure:tipo-estudios"

$ws.Range("F3").Value = "medida"

$ws.Range("D4").Value = "URI-Comunidad"
$ws.Range("F4").Value = "xsd:int"

# Remove row 5 entirely (mapping-aragon.xlsx / mapping-tipo-estudios.xlsx no longer needed)
$ws.Range("A5:G5").Delete()
